# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right after the "总计" (summary) sheet,
#    cloned from the "2022-Q1" template sheet so it inherits the same layout
#    and cell styles, then overwrite its data row with the new quarter's data.
# 2) Insert a new row into the "总计" summary sheet for "2022-Q3", pushing the
#    existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

$summary  = $wb.Worksheets.Item(1)   # "总计"
$template = $wb.Worksheets.Item(2)   # "2022-Q1" - used as a style/layout template

# --- Step 1: clone the template sheet right after "总计" and rename it ---
$template.Copy($null, $summary)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Overwrite the data row (row 2) with the new quarter's figures, keeping the
# same text-vs-number typing as the rest of the workbook: D/E/F/G are stored
# as text (to preserve formatting such as trailing zeros), H is numeric.
$c = $newSheet.Range("D2"); $c.NumberFormat = "@"; $c.Value = "1.12";   $c.ClearFormats()
$c = $newSheet.Range("E2"); $c.NumberFormat = "@"; $c.Value = "90.06";  $c.ClearFormats()
$c = $newSheet.Range("F2"); $c.NumberFormat = "@"; $c.Value = "3.57";   $c.ClearFormats()
$c = $newSheet.Range("G2"); $c.NumberFormat = "@"; $c.Value = "0.0400"; $c.ClearFormats()
$newSheet.Range("H2").Value = 2

# --- Step 2: insert a matching row into the "总计" summary sheet ---
$summary.Rows.Item(2).Insert()

# Give the new A2 the same style as the index column elsewhere (copy format
# from the row beneath, which just got pushed down from the old row 2).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$summary.Range("A2").Value = 0
$c = $summary.Range("B2"); $c.NumberFormat = "@";      $c.Value = "2022-Q3"; $c.ClearFormats()
$c = $summary.Range("C2"); $c.NumberFormat = "General"; $c.Value = 1;        $c.ClearFormats()
$c = $summary.Range("D2"); $c.NumberFormat = "General"; $c.Value = 0.04;     $c.ClearFormats()
